# Auto-generated script applying numeric corrections per the commit diff.
# Source: "chore: update Sheets via scheduled runner" - recalculated profit figures.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (11 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1813124.8
$ws.Range("J137").Value = 1820.421
$ws.Range("L137").Value = 5461.263
$ws.Range("N137").Value = -10561.263
$ws.Range("H141").Value = 2549.111
$ws.Range("I141").Value = 1230.2759
$ws.Range("J141").Value = 8012.857
$ws.Range("K141").Value = 3690.8277
$ws.Range("L141").Value = 24038.571
$ws.Range("M141").Value = 1489.1723
$ws.Range("N141").Value = -34398.571

# --- Sheet: ARM (46 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7153908
$ws.Range("I32").Value = 8073993
$ws.Range("K32").Value = 8073993
$ws.Range("M32").Value = -8073706
$ws.Range("H45").Value = 2821.0833
$ws.Range("I45").Value = 2800
$ws.Range("J45").Value = 2836.1428
$ws.Range("K45").Value = 2800
$ws.Range("L45").Value = 2836.1428
$ws.Range("M45").Value = -2423
$ws.Range("N45").Value = -3590.1428
$ws.Range("H61").Value = 10419960
$ws.Range("I61").Value = 16668656
$ws.Range("J61").Value = 5466.6665
$ws.Range("K61").Value = 16668656
$ws.Range("L61").Value = 5466.6665
$ws.Range("M61").Value = -16668444
$ws.Range("N61").Value = -5890.6665
$ws.Range("H74").Value = 33336524
$ws.Range("I74").Value = 2039.3334
$ws.Range("J74").Value = 55559516
$ws.Range("K74").Value = 2039.3334
$ws.Range("L74").Value = 55559516
$ws.Range("M74").Value = -1165.3334
$ws.Range("N74").Value = -55561264
$ws.Range("H77").Value = 33336524
$ws.Range("I77").Value = 2039.3334
$ws.Range("J77").Value = 55559516
$ws.Range("K77").Value = 10196.667
$ws.Range("L77").Value = 277797580
$ws.Range("M77").Value = -5828.666999999999
$ws.Range("N77").Value = -277806316
$ws.Range("H132").Value = 1085840
$ws.Range("I132").Value = 1837.0927
$ws.Range("J132").Value = 4529143
$ws.Range("K132").Value = 5511.2781
$ws.Range("L132").Value = 13587429
$ws.Range("M132").Value = -2981.2781
$ws.Range("N132").Value = -13592489
$ws.Range("H136").Value = 10419960
$ws.Range("I136").Value = 16668656
$ws.Range("J136").Value = 5466.6665
$ws.Range("K136").Value = 50005968
$ws.Range("L136").Value = 16399.9995
$ws.Range("M136").Value = -50003418
$ws.Range("N136").Value = -21499.9995

# --- Sheet: BSM (7 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2475.2354
$ws.Range("I134").Value = 2429.0303
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 7287.090899999999
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -4752.090899999999
$ws.Range("N134").Value = -17070

# --- Sheet: CRP (29 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 19428.857
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 25200.4
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 25200.4
$ws.Range("M4").Value = -4888
$ws.Range("N4").Value = -25424.4
$ws.Range("H31").Value = 6089.3237
$ws.Range("I31").Value = 1790.5
$ws.Range("J31").Value = 7339.891
$ws.Range("K31").Value = 1790.5
$ws.Range("L31").Value = 7339.891
$ws.Range("M31").Value = -1495.5
$ws.Range("N31").Value = -7929.891
$ws.Range("H34").Value = 6089.3237
$ws.Range("I34").Value = 1790.5
$ws.Range("J34").Value = 7339.891
$ws.Range("K34").Value = 1790.5
$ws.Range("L34").Value = 7339.891
$ws.Range("M34").Value = -1588.5
$ws.Range("N34").Value = -7743.891
$ws.Range("H58").Value = 1600
$ws.Range("I58").Value = 1600
$ws.Range("K58").Value = 1600
$ws.Range("M58").Value = -1397
$ws.Range("H136").Value = 1600
$ws.Range("I136").Value = 1600
$ws.Range("K136").Value = 4800
$ws.Range("M136").Value = -2250

# --- Sheet: CUL (40 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 975.2174
$ws.Range("I20").Value = 810
$ws.Range("K20").Value = 2430
$ws.Range("M20").Value = -2203
$ws.Range("H49").Value = 7786.6665
$ws.Range("J49").Value = 7786.6665
$ws.Range("L49").Value = 23359.9995
$ws.Range("N49").Value = -23671.9995
$ws.Range("H106").Value = 9823.727999999999
$ws.Range("J106").Value = 9823.727999999999
$ws.Range("L106").Value = 29471.184
$ws.Range("N106").Value = -31363.184
$ws.Range("H112").Value = 8947.833000000001
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 9670.362999999999
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 29011.089
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -31227.089
$ws.Range("H121").Value = 1045.7333
$ws.Range("I121").Value = 294.7
$ws.Range("J121").Value = 2547.8
$ws.Range("K121").Value = 884.0999999999999
$ws.Range("L121").Value = 7643.400000000001
$ws.Range("M121").Value = 425.9000000000001
$ws.Range("N121").Value = -10263.4
$ws.Range("H125").Value = 5250.6665
$ws.Range("I125").Value = 850
$ws.Range("J125").Value = 6130.8
$ws.Range("K125").Value = 2550
$ws.Range("L125").Value = 18392.4
$ws.Range("M125").Value = 2370
$ws.Range("N125").Value = -28232.4
$ws.Range("H129").Value = 778214.4
$ws.Range("I129").Value = 432.46155
$ws.Range("J129").Value = 1167105.2
$ws.Range("K129").Value = 1297.38465
$ws.Range("L129").Value = 3501315.6
$ws.Range("M129").Value = 3702.61535
$ws.Range("N129").Value = -3511315.6

# --- Sheet: GSM (8 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 47753.668
$ws.Range("J39").Value = 47753.668
$ws.Range("L39").Value = 47753.668
$ws.Range("N39").Value = -48817.668
$ws.Range("H109").Value = 26463.75
$ws.Range("J109").Value = 26463.75
$ws.Range("L109").Value = 26463.75
$ws.Range("N109").Value = -28543.75

# --- Sheet: LTW (22 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9085.714
$ws.Range("I2").Value = 5001
$ws.Range("J2").Value = 9766.5
$ws.Range("K2").Value = 5001
$ws.Range("L2").Value = 9766.5
$ws.Range("M2").Value = -4889
$ws.Range("N2").Value = -9990.5
$ws.Range("H22").Value = 7591.7334
$ws.Range("I22").Value = 780
$ws.Range("K22").Value = 780
$ws.Range("M22").Value = -485
$ws.Range("H27").Value = 7591.7334
$ws.Range("I27").Value = 780
$ws.Range("K27").Value = 780
$ws.Range("M27").Value = -673
$ws.Range("H46").Value = 550
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 537.5
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 537.5
$ws.Range("M46").Value = -412
$ws.Range("N46").Value = -913.5

# --- Sheet: WVR (4 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2739.1428
$ws.Range("I122").Value = 1929.1765
$ws.Range("K122").Value = 5787.529500000001
$ws.Range("M122").Value = -3337.529500000001
